$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 2. Data reporter / 1. Indicator information text updates
# ---------------------------------------------------------------------------

$ws.Range("B2").Value = '15. Protect, restore and promote sustainable use of terrestrial ecosystems, sustainably manage forests, combat desertification, and halt and reverse land degradation and halt biodiversity loss'

$ws.Range("B3").Value = '15.1 By 2020, ensure the conservation, restoration and sustainable use of terrestrial and inland freshwater ecosystems and their services, in particular forests, wetlands, mountains and drylands, in line with obligations under international agreements'

$ws.Range("B4").Value = '15.1.2 Proportion of important sites for terrestrial and freshwater biodiversity that are covered by protected areas, by ecosystem type'

$ws.Range("B6").Value = 'SAEPF, Department of Biodiversity Conservation and Specially Protected Natural Areas '

$ws.Range("B7").Value = 'Turdumatov Talantbek Kubanychbekovich, Musaev Almaz Mustafaevich'

$b8 = @'
info@fauna.kg,   
Turdumatov.fauna@gmail.com
'@
$ws.Range("B8").Value = $b8

$ws.Range("B9").Value = '+996 (312)   46-68-27'

$ws.Range("B10").Value = 'www.fauna.kg'

# B7 and B8 need wrap text turned on (previously unwrapped single-line style)
$ws.Range("B7").WrapText = $true
$ws.Range("B8").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Definitions and concepts
# ---------------------------------------------------------------------------

$ws.Range("B12").Value = 'Determines the share of land, water surface and airspace above them in the total territory of the republic, protected in accordance with national legislation, constituting the area of specially protected natural territories.'
$ws.Rows.Item(12).RowHeight = 45.75

$ws.Range("B13").Value = 'Specially protected natural territories are areas of territories and/or water areas having priority environmental, scientific, cultural, aesthetic and historical significance, which are national property, in whole or in part, permanently or temporarily withdrawn from economic activity, for which a special regime of protection and use is established.'
$ws.Rows.Item(13).RowHeight = 84.75

$ws.Range("B14").Value = 'They are created in order to preserve reference and unique natural complexes and objects, landmark natural formations, the genetic fund of the animal and plant world, study natural processes in the biosphere and control changes in its state.'
$ws.Rows.Item(14).RowHeight = 48

# ---------------------------------------------------------------------------
# 4. Data sources and collection methods
# ---------------------------------------------------------------------------

$b16 = @'
The data is generated in accordance with the resolutions of the Government of the Kyrgyz Republic and other regulatory legal acts on the formation of protected areas and their reorganization.
The data sources are departmental materials of the results of the work of specialized bodies on land and forestry (the State Agency for Land Resources under the PKR, the State Institution Kyrgyz Forest Management" under the PKR with the participation of the Institute of Biology of the National Academy of Sciences of the Kyrgyz Republic) on territories that are part of protected areas, as well as statistical reports of protected areas (state nature reserves and state natural parks) in the form No. 1-Protected areas (reports on specially protected natural areas)."
'@
$ws.Range("B16").Value = $b16
$ws.Rows.Item(16).RowHeight = 129

$ws.Range("B17").Value = 'The methods of data collection are the interdepartmental exchange of relevant information at the request of interested parties and statistical reporting.'
$ws.Rows.Item(17).RowHeight = 47.25

# ---------------------------------------------------------------------------
# 5. Method of computation and other methodological considerations
# ---------------------------------------------------------------------------

$b19 = @'
The percentage of protected areas from the total area of the country. 
The share of the area of protected areas in the total area of the country is the specific weight of the total area of protected areas in relation to the total area of the country.
The method of calculating this form of representation of the indicator:
Share of protected areas (%) = Total area of protected areas in ha / Area of the country in ha * 100

'@
$ws.Range("B19").Value = $b19
$ws.Rows.Item(19).RowHeight = 81

$ws.Range("B20").Value = 'There are likely discrepancies in indicators due to inconsistencies between data on the area of a particular protected area object, according to the regulatory legal act characterizing it, and data based on subsequent processing of cartographic materials for the corresponding object.'
$ws.Rows.Item(20).RowHeight = 84.75

$ws.Range("B21").Value = 'Timely introduction of amendments and additions to the regulatory legal acts on protected areas in relation to the clarification of their boundaries and areas.'

# ---------------------------------------------------------------------------
# 6. Data availability and disaggregation
# ---------------------------------------------------------------------------

$ws.Range("B23").Value = 'The data is available to interested parties on the websites of specialized bodies and official statistics, as well as country reports on biodiversity.'

# ---------------------------------------------------------------------------
# 8. References and documentation
# ---------------------------------------------------------------------------

$ws.Range("B26").Value = 'The indicator is formed according to the UNECE Guidelines "Environmental Indicators and Assessment Reports based on them", and the methodology of the SDG indicators, based on available data and consultations.'
$ws.Rows.Item(26).RowHeight = 65.25

# ---------------------------------------------------------------------------
# Selection moves to B24 (last active cell recorded in the file)
# ---------------------------------------------------------------------------
$ws.Range("B24").Select()
